# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Sat Mar  9 05:40:03 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '68.225.23'
$ws.Range('E2').Value = "'" + '  +2.06%  '
$ws.Range('D3').Value = "'" + '3.922.42'
$ws.Range('E3').Value = "'" + '  +1.06%  '
$ws.Range('E4').Value = "'" + '  -0.05%  '
$ws.Range('D5').Value = "'" + '488.59'
$ws.Range('E5').Value = "'" + '  +4.30%  '
$ws.Range('D6').Value = "'" + '147.43'
$ws.Range('E6').Value = "'" + '  +3.01%  '
$ws.Range('E7').Value = "'" + '  +0.94%  '
$ws.Range('E8').Value = "'" + '  -0.03%  '
$ws.Range('E9').Value = "'" + '  +0.16%  '
$ws.Range('E10').Value = "'" + '  +6.05%  '
$ws.Range('D11').Value = "'" + '0.0000356'
$ws.Range('E11').Value = "'" + '  +8.20%  '
$ws.Range('D12').Value = "'" + '42.75'
$ws.Range('E12').Value = "'" + '  -0.04%  '
$ws.Range('D13').Value = "'" + '10.62'
$ws.Range('E13').Value = "'" + '  +3.17%  '
$ws.Range('D14').Value = "'" + '4.545.69'
$ws.Range('E14').Value = "'" + '  +0.83%  '
$ws.Range('D15').Value = "'" + '14.78'
$ws.Range('E15').Value = "'" + '  -0.58%  '
$ws.Range('D16').Value = "'" + '3.921.78'
$ws.Range('E16').Value = "'" + '  +1.33%  '
$ws.Range('E17').Value = "'" + '  -0.05%  '
$ws.Range('D18').Value = "'" + '20.01'
$ws.Range('E18').Value = "'" + '  +1.24%  '
$ws.Range('E19').Value = "'" + '  -1.37%  '
$ws.Range('D20').Value = "'" + '68.334.61'
$ws.Range('D21').Value = "'" + '444.15'
$ws.Range('E21').Value = "'" + '  +3.97%  '
$ws.Range('D22').Value = "'" + '14.74'
$ws.Range('E22').Value = "'" + '  +0.63%  '
$ws.Range('E23').Value = "'" + '  +1.94%  '
$ws.Range('D24').Value = "'" + '88.51'
$ws.Range('E24').Value = "'" + '  +0.75%  '
$ws.Range('D25').Value = "'" + '11.67'
$ws.Range('E25').Value = "'" + '  +17.42%  '
$ws.Range('D26').Value = "'" + '11.28'
$ws.Range('E26').Value = "'" + '  +17.40%  '
$ws.Range('D27').Value = "'" + '3.63'
$ws.Range('E27').Value = "'" + '  +2.55%  '
$ws.Range('D28').Value = "'" + '38.90'
$ws.Range('E28').Value = "'" + '  +1.83%  '
$ws.Range('D29').Value = "'" + '5.88'
$ws.Range('E29').Value = "'" + '  +1.99%  '
$ws.Range('D30').Value = "'" + '715.28'
$ws.Range('E30').Value = "'" + '  -1.82%  '
$ws.Range('E31').Value = "'" + '  -1.55%  '
$ws.Range('E32').Value = "'" + '  +1.10%  '
$ws.Range('E33').Value = "'" + '  +2.58%  '
$ws.Range('D34').Value = "'" + '0.0₃0909'
$ws.Range('E34').Value = "'" + '  +18.78%  '
$ws.Range('D35').Value = "'" + '41.33'
$ws.Range('E35').Value = "'" + '  -4.20%  '
$ws.Range('D36').Value = "'" + '5.89'
$ws.Range('E36').Value = "'" + '  +9.50%  '
$ws.Range('D37').Value = "'" + '59.05'
$ws.Range('E37').Value = "'" + '  +2.88%  '
$ws.Range('E38').Value = "'" + '  -4.43%  '
$ws.Range('E39').Value = "'" + '  +0.23%  '
$ws.Range('B40').Value = "'" + 'TheGraph'
$ws.Range('C40').Value = "'" + 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').Value = "'" + '0.389'
$ws.Range('E40').Value = "'" + '  +15.98%  '
$ws.Range('B41').Value = "'" + 'Fetch.AI'
$ws.Range('C41').Value = "'" + 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').Value = "'" + '2.92'
$ws.Range('E41').Value = "'" + '  +14.98%  '
$ws.Range('E42').Value = "'" + '  +1.48%  '
$ws.Range('D43').Value = "'" + '3.13'
$ws.Range('E43').Value = "'" + '  +1.68%  '
$ws.Range('D44').Value = "'" + '2.92'
$ws.Range('E44').Value = "'" + '  +4.60%  '
$ws.Range('E45').Value = "'" + '  +1.94%  '
$ws.Range('E46').Value = "'" + '  -0.09%  '
$ws.Range('D47').Value = "'" + '3.42'
$ws.Range('E47').Value = "'" + '  +1.23%  '
$ws.Range('D48').Value = "'" + '2.14'
$ws.Range('E48').Value = "'" + '  -0.45%  '
$ws.Range('B49').Value = "'" + 'BabyDogeCoin'
$ws.Range('C49').Value = "'" + 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = "'" + '0.0₆0343'
$ws.Range('E49').Value = "'" + '  +44.93%  '
$ws.Range('B50').Value = "'" + 'Monero'
$ws.Range('C50').Value = "'" + 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = "'" + '145.11'
$ws.Range('E50').Value = "'" + '  +1.07%  '
$ws.Range('D51').Value = "'" + '3.12'
$ws.Range('E51').Value = "'" + '  +0.32%  '
